# Update the Macro_taxonomy mapping table: split the "Other" category
# (for both Urban and Rural settlements) into two taxonomy rows, each
# carrying half of the original proportion, then leave the UI selection
# state matching the author's last interaction (Macro_taxonomy tab active
# with D33 selected; Costs sheet selection left at B24).

$wb = $excel.ActiveWorkbook

$wsMacro = $wb.Worksheets.Item("Macro_taxonomy")
$wsCosts = $wb.Worksheets.Item("Costs")

# --- Urban / "Other" split (row 17 -> rows 17 & 18) ---------------------
$wsMacro.Rows.Item(18).Insert()

$wsMacro.Range("D17").Value = 0.5

$wsMacro.Range("A18").Value = "Other"
$wsMacro.Range("B18").Value = "Urban"
$wsMacro.Range("C18").Value = "ME+MEO/LWAL"
$wsMacro.Range("D18").Value = 0.5

# --- Rural / "Other" split (old row 30, now row 31 -> rows 31 & 32) -----
$wsMacro.Rows.Item(32).Insert()

$wsMacro.Range("D31").Value = 0.5

$wsMacro.Range("A32").Value = "Other"
$wsMacro.Range("B32").Value = "Rural"
$wsMacro.Range("C32").Value = "ME+MEO/LWAL"
$wsMacro.Range("D32").Value = 0.5

# --- Leave the selection / active-sheet state as the author left it -----
$wsCosts.Range("B24").Select()

$wsMacro.Activate()
$wsMacro.Range("D33").Select()
